# Rock_Stars_Presentation.pptx - "small fix in presentation2"
#
# The title text box on slide 1 currently reads:
#   "SAFE TRAVEL " + "POPULATION " + "HEALTH DOMAIN"
# (first run "SAFE TRAVEL ", line break, then "POPULATION HEALTH DOMAIN").
#
# The fix corrects "TRAVEL" -> "TRAVELS", splitting the original first run
# into two runs: "SAFE " and "TRAVELS ".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# Characters 6-12 of the title text are "TRAVEL " (the word right after
# "SAFE "). Replacing just that sub-range keeps "SAFE " as its own run
# and creates a new run for the replacement text, turning
# "SAFE TRAVEL " into "SAFE " + "TRAVELS ".
$tr.Characters(6, 7).Text = "TRAVELS "
